# Updated cryptos list on Tue Feb 20 20:40:32 UTC 2024 with GitHub Actions
# Refresh the Price (col D) and Volume(1h) (col E) cells on the active sheet
# to the latest scraped values. Price cells whose new text would otherwise be
# auto-parsed by Excel as a number (losing formatting, e.g. "1.00" -> 1) are
# forced to a Text number format first so the literal string survives.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "52.156.19"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3
$ws.Range("D3").Value = "2.975.00"
$ws.Range("E3").Value = "  +1.36%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.38"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.19"
$ws.Range("E6").Value = "  -4.66%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.614"
$ws.Range("E9").Value = "  -1.94%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.28"
$ws.Range("E10").Value = "  -2.89%  "

# Row 11
$ws.Range("E11").Value = "  +1.33%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0857"
$ws.Range("E12").Value = "  -3.76%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.22"
$ws.Range("E13").Value = "  -3.91%  "

# Row 14
$ws.Range("D14").Value = "3.439.46"
$ws.Range("E14").Value = "  +1.19%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.69"
$ws.Range("E15").Value = "  -1.72%  "

# Row 16
$ws.Range("D16").Value = "2.969.86"
$ws.Range("E16").Value = "  +1.22%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.00"
$ws.Range("E17").Value = "  +1.40%  "

# Row 18
$ws.Range("D18").Value = "52.166.90"
$ws.Range("E18").Value = "  +0.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.49"
$ws.Range("E19").Value = "  +5.24%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.48"
$ws.Range("E20").Value = "  -2.11%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.60"
$ws.Range("E21").Value = "  -4.69%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0974"
$ws.Range("E22").Value = "  -1.26%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.58"
$ws.Range("E23").Value = "  -2.32%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.18"
$ws.Range("E24").Value = "  -1.81%  "

# Row 25
$ws.Range("E25").Value = "  -1.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.179"
$ws.Range("E26").Value = "  +0.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.82"
$ws.Range("E27").Value = "  -0.53%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.55"
$ws.Range("E28").Value = "  +0.99%  "

# Row 29
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("E30").Value = "  +0.82%  "

# Row 31
$ws.Range("E31").Value = "  -2.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.12"
$ws.Range("E32").Value = "  -1.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.19"
$ws.Range("E33").Value = "  -3.35%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.17"
$ws.Range("E34").Value = "  -4.31%  "

# Row 35
$ws.Range("E35").Value = "  -3.98%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0443"
$ws.Range("E36").Value = "  -2.50%  "

# Row 37
$ws.Range("E37").Value = "  -0.05%  "

# Row 38
$ws.Range("E38").Value = "  -3.23%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.91"
$ws.Range("E39").Value = "  -5.31%  "

# Row 40
$ws.Range("E40").Value = "  -4.13%  "

# Row 41
$ws.Range("E41").Value = "  +0.32%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  -0.78%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.72"
$ws.Range("E43").Value = "  -1.82%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.77"
$ws.Range("E44").Value = "  +8.99%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.13"
$ws.Range("E45").Value = "  -3.18%  "

# Row 46
$ws.Range("D46").Value = "2.117.22"
$ws.Range("E46").Value = "  -2.47%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.39"
$ws.Range("E47").Value = "  -4.10%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  -7.45%  "

# Row 49
$ws.Range("E49").Value = "  -3.47%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0337"
$ws.Range("E50").Value = "  -2.67%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.935"
$ws.Range("E51").Value = "  -0.47%  "
